$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B15").Value = "write php script to reset everything with app"
$ws.Range("B15").Select()
